$d = $word.ActiveDocument

$replacements = @(
    @("865÷3=", "215÷5="),
    @("929÷4=", "632÷5="),
    @("946÷6=", "889÷7="),
    @("330÷8=", "376÷8="),
    @("535÷6=", "120÷5="),
    @("683÷9=", "899÷3="),
    @("945÷8=", "440÷3="),
    @("122÷3=", "844÷8="),
    @("811÷8=", "389÷4="),
    @("691÷5=", "361÷6="),
    @("112÷3=", "255÷6="),
    @("114÷6=", "620÷3="),
    @("959÷4=", "497÷9="),
    @("472÷9=", "628÷7="),
    @("458÷7=", "236÷2="),
    @("589÷9=", "865÷4="),
    @("845÷3=", "823÷3="),
    @("603÷6=", "166÷6="),
    @("410÷9=", "525÷6="),
    @("974÷4=", "676÷3="),
    @("808÷9=", "807÷6="),
    @("697÷4=", "997÷6="),
    @("125÷8=", "804÷4="),
    @("138÷3=", "262÷6="),
    @("339÷5=", "211÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
